# Cleans up the "markdown residue" left by the original authoring tool:
# headers like "**Study Overview:*"+"*" (split across runs by Word's
# grammar-checker proofErr markers) collapse to a single run reading
# "**Study Overview:**", and quoted "PaGamO" mentions (split by
# spell-check proofErr markers) collapse to a single run as well.
# Finally, the trailing "This detailed summary..." paragraph is cleared
# to an empty paragraph.

$d = $word.ActiveDocument
$placeholder = [char]0x01
$q = [char]34

function Set-ParagraphText {
    param(
        [int]$Index,
        [string]$NewText
    )

    $para = $d.Paragraphs.Item($Index)
    $rng = $para.Range
    # Exclude the trailing paragraph mark from the range.
    $rng.End = $rng.End - 1

    if ($rng.Start -eq $rng.End -and $NewText -eq "") {
        return
    }

    # First pass collapses every run (and any proofErr markers) in the
    # paragraph down to a single run holding a 1-character placeholder.
    $rng.Text = $placeholder

    # Second pass (re-fetch the range - the paragraph now has a single
    # run) sets the final text on that lone run.
    $rng2 = $para.Range
    $rng2.End = $rng2.End - 1
    $rng2.Text = $NewText
}

Set-ParagraphText 2 "**Study Overview:**"
$t3 = "- The study focuses on the effectiveness of an educational game called ${q}PaGamO${q} in enhancing student learning in higher education."
Set-ParagraphText 3 $t3

Set-ParagraphText 7 "**Intrinsic Motivation:**"

Set-ParagraphText 12 "**Objective:**"
$t13 = "- The study's objective is to assess whether the use of the educational game ${q}PaGamO${q} could enhance students' learning ability and understand students' perceptions of educational games."
Set-ParagraphText 13 $t13

Set-ParagraphText 15 "**Materials and Methods:**"
$t17 = "- The ${q}PaGamO${q} game was introduced as a supplementary tool for learning."
Set-ParagraphText 17 $t17

Set-ParagraphText 21 "**Data Collection and Analysis:**"
$t22 = "- The study used a mixed-method approach to assess the effectiveness of ${q}PaGamO.${q}"
Set-ParagraphText 22 $t22
$t25 = "- The analysis included factors such as students' ${q}PaGamO${q} scores, examination scores, motives for playing, and perceptions of the game's effectiveness."
Set-ParagraphText 25 $t25

Set-ParagraphText 27 "**Results and Discussion:**"
$t28 = "- The study found a significant relationship between ${q}PaGamO${q} scores and multiple-choice (MC) exam scores."
Set-ParagraphText 28 $t28
$t29 = "- Students' motives for playing ${q}PaGamO${q} were primarily intrinsic, with factors like fun, self-learning, and wanting to perform well in the final examination being prominent."
Set-ParagraphText 29 $t29
$t32 = "- The choice of devices for playing ${q}PaGamO${q} varied, with mobile phones being the most popular due to their convenience."
Set-ParagraphText 32 $t32
$t33 = "- The study showed that even short sessions of using ${q}PaGamO${q} were effective in helping students prepare for exams."
Set-ParagraphText 33 $t33

Set-ParagraphText 35 "**Recommendations for Further Research:**"

Set-ParagraphText 40 "**Conclusion:**"

Set-ParagraphText 45 "**Funding and Ethics:**"

# Clear the trailing "This detailed summary..." paragraph down to an
# empty paragraph (matches the now-empty <w:p/> before the sectPr).
$lastIndex = $d.Paragraphs.Count
Set-ParagraphText $lastIndex ""

Write-Output "done"
